$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Type" classification column (E) replacing the old "Row"/"Column"
# reference columns (E/F). Column F's old values are cleared (kept blank,
# formatted) while column E gets re-labelled for every data row.
$typeValues = @(
    "Type",           # row 1  (header)
    "Primary",        # row 2
    "Primary",        # row 3
    "Secondary",      # row 4
    "Secondary",      # row 5
    "Secondary",      # row 6
    "Secondary",      # row 7
    "Finished_Sum",   # row 8
    "Finished_Sum",   # row 9
    "Finished_Sum",   # row 10
    "Finished_Sum",   # row 11
    "Finished",       # row 12
    "Secondary_Sum",  # row 13
    "Finished",       # row 14
    "Finished",       # row 15
    "Finished",       # row 16
    "Finished",       # row 17
    "Finished",       # row 18
    "Finished",       # row 19
    "Finished",       # row 20
    "Finished",       # row 21
    "Finished",       # row 22
    "Finished",       # row 23
    "Finished",       # row 24
    "Finished",       # row 25
    "Finished",       # row 26
    "Finished",       # row 27
    "Finished"        # row 28
)

for ($i = 0; $i -lt $typeValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = $typeValues[$i]
}

# Rows 1-13 already carried the "Text" number format (style index 1) on
# column E; rows 14-28 are brand-new cells in column E and need the same
# "Text" format applied so they match the rest of the column.
$ws.Range("E1:E28").NumberFormat = "@"

# Column F ("Column") is no longer used for rows 1-13; clear its contents
# but leave the cell formatting (style) intact.
$ws.Range("F1:F13").ClearContents()

# Move the selection to match the saved view state (E22).
$ws.Range("E22").Select()
